$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-16 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-17 Sunday", 2) | Out-Null
$d.Content.Find.Execute("541÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "781÷3=", 2) | Out-Null
$d.Content.Find.Execute("784÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "737÷7=", 2) | Out-Null
$d.Content.Find.Execute("660÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "663÷9=", 2) | Out-Null
$d.Content.Find.Execute("664÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "441÷9=", 2) | Out-Null
$d.Content.Find.Execute("101÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "469÷7=", 2) | Out-Null
$d.Content.Find.Execute("781÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "708÷7=", 2) | Out-Null
$d.Content.Find.Execute("517÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "862÷7=", 2) | Out-Null
$d.Content.Find.Execute("352÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "173÷6=", 2) | Out-Null
$d.Content.Find.Execute("594÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "271÷8=", 2) | Out-Null
$d.Content.Find.Execute("587÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "567÷7=", 2) | Out-Null
$d.Content.Find.Execute("923÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "360÷3=", 2) | Out-Null
$d.Content.Find.Execute("211÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "835÷9=", 2) | Out-Null
$d.Content.Find.Execute("183÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "833÷9=", 2) | Out-Null
$d.Content.Find.Execute("923÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "749÷8=", 2) | Out-Null
$d.Content.Find.Execute("556÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "105÷9=", 2) | Out-Null
$d.Content.Find.Execute("208÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "507÷6=", 2) | Out-Null
$d.Content.Find.Execute("194÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "384÷2=", 2) | Out-Null
$d.Content.Find.Execute("740÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "554÷4=", 2) | Out-Null
$d.Content.Find.Execute("381÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "860÷4=", 2) | Out-Null
$d.Content.Find.Execute("903÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "484÷4=", 2) | Out-Null
$d.Content.Find.Execute("459÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "794÷6=", 2) | Out-Null
$d.Content.Find.Execute("269÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "716÷6=", 2) | Out-Null
$d.Content.Find.Execute("196÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "712÷3=", 2) | Out-Null
$d.Content.Find.Execute("159÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "769÷3=", 2) | Out-Null
$d.Content.Find.Execute("336÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "179÷8=", 2) | Out-Null
